$d = $word.ActiveDocument

# Helper: replace the text of a whole paragraph (minus its trailing paragraph
# mark) using InsertXML. Unlike Range.Text / Find.Execute, InsertXML only
# overwrites the exact target Range, so any *other* runs in the paragraph
# (e.g. the leading empty <w:r/> runs used throughout this document) are left
# completely untouched. $rPrXml lets callers preserve run formatting
# (bold/italic/etc.) on the replaced run.
function Set-ParagraphText {
    param(
        [int]$Index,
        [string]$OldText,
        [string]$NewText,
        [string]$RPrXml = ""
    )
    $p = $d.Paragraphs.Item($Index)
    $startPos = $p.Range.Start
    $endPos = $p.Range.End - 1   # exclude the paragraph mark
    $rr = $d.Range($startPos, $endPos)
    if ($rr.Text -ne $OldText) {
        throw "Set-ParagraphText: paragraph $Index text mismatch. Expected [$OldText] but found [$($rr.Text)]"
    }
    $escaped = $NewText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $tTag = "<w:t>" + $escaped + "</w:t>"
    if ($NewText.StartsWith(" ") -or $NewText.EndsWith(" ")) {
        $tTag = '<w:t xml:space="preserve">' + $escaped + "</w:t>"
    }
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p><w:r>' + $RPrXml + $tTag + '</w:r></w:p></w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $rr.InsertXML($xml)
}

# 1. Title (Heading1) - paragraph 1 (single run, no surrounding empty runs)
$titlePara = $d.Paragraphs.Item(1)
$titleCheckRange = $d.Range($titlePara.Range.Start, $titlePara.Range.End - 1)
if ($titleCheckRange.Text -ne "Play Candyways Bonanza 2 Megaways Free: Slot Game Review") {
    throw "Paragraph 1 text mismatch: [$($titleCheckRange.Text)]"
}
$titlePara.Range.Find.Execute(
    "Play Candyways Bonanza 2 Megaways Free: Slot Game Review", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Candyways Bonanza 2 Megaways for Free", 2) | Out-Null

# 2. "What we like" bullet list (each paragraph has a leading empty <w:r/>)
# NOTE: this runtime's PowerShell subset does not support named (-Param value)
# argument binding reliably, so all calls below use positional arguments.
Set-ParagraphText 39 "Variety of symbols and bonuses" "Variety of symbols and extra features"
Set-ParagraphText 40 "117,649 winning combinations" "Impressive number of winning combinations"
Set-ParagraphText 41 "Affordable betting range" "Wide betting range"
Set-ParagraphText 42 "Playable on both PC and mobile" "Compatible on PC and mobile devices"

# 3. "What we don't like" bullet list
Set-ParagraphText 44 "Limited betting range" "May not appeal to players who prefer traditional slot game themes"
Set-ParagraphText 45 "No progressive jackpot feature" "Free spins feature may be difficult to trigger"

# 4. Bold closing headline (paragraph 46) - preserve bold formatting
Set-ParagraphText 46 "Play Candyways Bonanza 2 Megaways Free: Slot Game Review" "Play Candyways Bonanza 2 Megaways for Free" "<w:rPr><w:b/></w:rPr>"

# 5. Italic meta description (paragraph 47) - preserve italic formatting
Set-ParagraphText 47 "Get the lowdown on Candyways Bonanza 2 Megaways slot game with our review. Play for free and experience its variety of symbols, winning combos, and bonuses." "Read our review of Candyways Bonanza 2 Megaways and play this exciting slot game for free." "<w:rPr><w:i/></w:rPr>"
